# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) values on the zh-cn and de-de
# sheets to reflect a newer report-generation run.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-25 01:30:55"
$zhcn.Range("H2").Value = "2016-03-25 01:31:23"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-25 01:31:00"
$dede.Range("H2").Value = "2016-03-25 01:31:34"
